$d = $word.ActiveDocument

# 1. "Alstom (Cetic)" -> "Cetic" (unique occurrence in the attendee table)
$d.Content.Find.Execute("Alstom (Cetic)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cetic", 2) | Out-Null

# 2. Second "NS" occurrence (the one right after "Jan Welvaarts") -> "L'loyds Register Rail"
#    There are two "NS" cells in the attendee table; anchor on "Jan Welvaarts" first,
#    then search forward from there so only the second "NS" is touched.
$rng = $d.Content
$rng.Find.Execute("Jan Welvaarts", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.Find.Execute("NS", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($rng.Start, $rng.End)
$target.Text = "L’loyds Register Rail"

# 3. Footer page-number field result: cached "1" -> "2" (the PAGE field in the
#    default/primary footer), keeping the field codes intact - only the cached
#    literal text inside the field result run is updated.
$ftr = $d.Sections(1).Footers(1)
for ($i = 1; $i -le $ftr.Range.Fields.Count; $i++) {
    $fld = $ftr.Range.Fields($i)
    if ($fld.Code.Text.Trim() -eq "PAGE") {
        $pageResult = $fld.Result
        $pageResult.Find.Execute("1", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "2", 2) | Out-Null
        break
    }
}
